$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '27.354.26'
Set-TextValue 'E2' '  -0.66%  '
Set-TextValue 'D3' '1.708.95'
Set-TextValue 'E3' '  -0.98%  '
Set-TextValue 'E4' '  -0.12%  '
Set-TextValue 'D5' '224.15'
Set-TextValue 'E5' '  -0.62%  '
Set-TextValue 'D6' '0.5296'
Set-TextValue 'E6' '  -1.34%  '
Set-TextValue 'E7' '  -0.06%  '
Set-TextValue 'E8' '  -0.33%  '
Set-TextValue 'D9' '0.06618'
Set-TextValue 'E9' '  +0.22%  '
Set-TextValue 'D10' '20.82'
Set-TextValue 'E10' '  -4.42%  '
Set-TextValue 'D11' '0.07664'
Set-TextValue 'E11' '  -0.64%  '
Set-TextValue 'D12' '4.509'
Set-TextValue 'E12' '  -2.19%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.729.68'
Set-TextValue 'E13' '  +0.04%  '
Set-TextValue 'B14' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C14' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D14' '1.943.06'
Set-TextValue 'E14' '  -1.08%  '
Set-TextValue 'D15' '0.5818'
Set-TextValue 'E15' '  -0.49%  '
Set-TextValue 'D16' '0.0₅8172'
Set-TextValue 'E16' '  -1.59%  '
Set-TextValue 'D17' '67.79'
Set-TextValue 'E17' '  -0.34%  '
Set-TextValue 'D18' '27.348.28'
Set-TextValue 'E18' '  -0.77%  '
Set-TextValue 'D19' '217.64'
Set-TextValue 'E19' '  -1.77%  '
Set-TextValue 'E20' '  -0.05%  '
Set-TextValue 'D21' '4.631'
Set-TextValue 'E21' '  -2.18%  '
Set-TextValue 'E22' '  -2.19%  '
Set-TextValue 'D23' '5.980'
Set-TextValue 'E23' '  -1.89%  '
Set-TextValue 'D24' '1.004'
Set-TextValue 'E24' '  -0.10%  '
Set-TextValue 'D25' '143.88'
Set-TextValue 'E25' '  -3.05%  '
Set-TextValue 'D26' '1.691'
Set-TextValue 'E26' '  -1.31%  '
Set-TextValue 'E27' '  -2.34%  '
Set-TextValue 'D28' '7.259'
Set-TextValue 'E28' '  -1.98%  '
Set-TextValue 'D29' '16.23'
Set-TextValue 'E29' '  -2.62%  '
Set-TextValue 'D30' '0.05370'
Set-TextValue 'E30' '  -3.66%  '
Set-TextValue 'D31' '1.293'
Set-TextValue 'E31' '  -0.60%  '
Set-TextValue 'D32' '3.468'
Set-TextValue 'E32' '  -2.35%  '
Set-TextValue 'D33' '3.422'
Set-TextValue 'E33' '  -0.99%  '
Set-TextValue 'D34' '1.647'
Set-TextValue 'E34' '  -0.77%  '
Set-TextValue 'D35' '2.868'
Set-TextValue 'E35' '  +1.40%  '
Set-TextValue 'D36' '0.9514'
Set-TextValue 'E36' '  -1.27%  '
Set-TextValue 'D37' '2.397'
Set-TextValue 'E37' '  -2.00%  '
Set-TextValue 'D38' '0.5861'
Set-TextValue 'E38' '  -1.41%  '
Set-TextValue 'D39' '0.01639'
Set-TextValue 'E39' '  -0.52%  '
Set-TextValue 'B40' 'Maker'
Set-TextValue 'C40' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D40' '1.062.13'
Set-TextValue 'E40' '  +0.65%  '
Set-TextValue 'B41' 'FraxShare'
Set-TextValue 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D41' '5.803'
Set-TextValue 'E41' '  -2.14%  '
Set-TextValue 'D42' '0.8456'
Set-TextValue 'E42' '  -1.53%  '
Set-TextValue 'D43' '1.004'
Set-TextValue 'E43' '  -0.04%  '
Set-TextValue 'D44' '100.88'
Set-TextValue 'E44' '  -0.59%  '
Set-TextValue 'D45' '1.851.18'
Set-TextValue 'E45' '  -0.99%  '
Set-TextValue 'D46' '0.0₈8114'
Set-TextValue 'E46' '  -0.41%  '
Set-TextValue 'D47' '57.76'
Set-TextValue 'E47' '  -2.25%  '
Set-TextValue 'E48' '  +1.90%  '
Set-TextValue 'D49' '1.004'
Set-TextValue 'E49' '  +0.30%  '
Set-TextValue 'D50' '8.112'
Set-TextValue 'E50' '  -1.43%  '
Set-TextValue 'D51' '0.05231'
Set-TextValue 'E51' '  -0.61%  '
